# Competition.xlsx — add new competitor row "Asana" to the Table2 listing
# on Sheet1, and resize column I to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow the existing table (Table2, A1:J33 -> A1:J34) by one row so the
# table/autoFilter ranges and the sheet dimension all expand together.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$row = 34

# Columns: A Competition, B Category, C Website, D State, E HQ,
#          F Investors, G Management, H Platforms, I Business Model, J Notes
# Write G/I/J/H last (in that order) so new shared-string entries are
# appended in the same order the source workbook uses.
$ws.Cells.Item($row, 1).Value = "Asana"
$ws.Cells.Item($row, 2).Value = "Task management"
$ws.Cells.Item($row, 3).Value = "http://www.asana.com/"
$ws.Cells.Item($row, 4).Value = "GA"
$ws.Cells.Item($row, 7).Value = "Dustin Moskowitz"
$ws.Cells.Item($row, 9).Value = "Freemium? (30 free)"
$ws.Cells.Item($row, 10).Value = "slick web app - for teams / business"
$ws.Cells.Item($row, 8).Value = "Web, mobile"

# Column I needs to be a bit wider to fit "Freemium? (30 free)".
$ws.Columns.Item(9).ColumnWidth = 19.2

# Leave the selection where the author's commit left it.
$ws.Range("H35").Select()
